$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M2").Value = -269.7143
$ws.Range("H2").Value = 684.2727
$ws.Range("K2").Value = 382.7143
$ws.Range("I2").Value = 382.7143
$ws.Range("I11").Value = 2514
$ws.Range("H11").Value = 2514
$ws.Range("K11").Value = 2514
$ws.Range("M11").Value = -2374
$ws.Range("H16").Value = 91000
$ws.Range("K16").Value = 91000
$ws.Range("I16").Value = 91000
$ws.Range("M16").Value = -90770
$ws.Range("N17").Value = -16232.619
$ws.Range("H17").Value = 5246.0415
$ws.Range("L17").Value = 15896.619
$ws.Range("J17").Value = 5298.873
$ws.Range("H29").Value = 1444.6
$ws.Range("K29").Value = 4333.799999999999
$ws.Range("I29").Value = 1444.6
$ws.Range("M29").Value = -4052.799999999999
$ws.Range("K33").Value = 1186
$ws.Range("M33").Value = -957
$ws.Range("I33").Value = 1186
$ws.Range("H33").Value = 1125.5555
$ws.Range("N38").Value = -5332.5
$ws.Range("H38").Value = 288.07693
$ws.Range("L38").Value = 4588.5
$ws.Range("J38").Value = 1529.5
$ws.Range("J41").Value = 1928.5714
$ws.Range("K41").Value = 1875.25
$ws.Range("N41").Value = -2808.5714
$ws.Range("I41").Value = 1875.25
$ws.Range("M41").Value = -1435.25
$ws.Range("H41").Value = 1909.1818
$ws.Range("L41").Value = 1928.5714
$ws.Range("K51").Value = 6600
$ws.Range("I51").Value = 6600
$ws.Range("M51").Value = -6116
$ws.Range("H51").Value = 7246.5
$ws.Range("I53").Value = 340.5
$ws.Range("H53").Value = 7495.643
$ws.Range("K53").Value = 340.5
$ws.Range("M53").Value = 296.5
$ws.Range("K64").Value = 3975
$ws.Range("N64").Value = -8796.200000000001
$ws.Range("I64").Value = 3975
$ws.Range("M64").Value = -3727
$ws.Range("H64").Value = 7579.3335
$ws.Range("L64").Value = 8300.200000000001
$ws.Range("J64").Value = 8300.200000000001
$ws.Range("I67").Value = 3975
$ws.Range("N67").Value = -10016.2
$ws.Range("M67").Value = -3117
$ws.Range("H67").Value = 7579.3335
$ws.Range("L67").Value = 8300.200000000001
$ws.Range("J67").Value = 8300.200000000001
$ws.Range("K67").Value = 3975
$ws.Range("J70").Value = 5497
$ws.Range("K70").Value = 12225.6666
$ws.Range("N70").Value = -17031
$ws.Range("I70").Value = 4075.2222
$ws.Range("M70").Value = -11955.6666
$ws.Range("H70").Value = 4333.727
$ws.Range("L70").Value = 16491
$ws.Range("K73").Value = 12225.6666
$ws.Range("N73").Value = -18363
$ws.Range("I73").Value = 4075.2222
$ws.Range("M73").Value = -11289.6666
$ws.Range("H73").Value = 4333.727
$ws.Range("L73").Value = 16491
$ws.Range("J73").Value = 5497
$ws.Range("H74").Value = 7649.2
$ws.Range("L74").Value = 7857.826
$ws.Range("J74").Value = 7857.826
$ws.Range("N74").Value = -9729.826000000001
$ws.Range("I76").Value = 5918.154
$ws.Range("M76").Value = -5603.154
$ws.Range("H76").Value = 5891.1
$ws.Range("L76").Value = 5840.857
$ws.Range("J76").Value = 5840.857
$ws.Range("K76").Value = 5918.154
$ws.Range("N76").Value = -6470.857
$ws.Range("N77").Value = -48649.13
$ws.Range("H77").Value = 7649.2
$ws.Range("L77").Value = 39289.13
$ws.Range("J77").Value = 7857.826
$ws.Range("J79").Value = 5840.857
$ws.Range("K79").Value = 5918.154
$ws.Range("N79").Value = -8024.857
$ws.Range("I79").Value = 5918.154
$ws.Range("M79").Value = -4826.154
$ws.Range("H79").Value = 5891.1
$ws.Range("L79").Value = 5840.857
$ws.Range("N86").Value = -4447
$ws.Range("I86").Value = 2215.6667
$ws.Range("M86").Value = -1092.6667
$ws.Range("L86").Value = 2201
$ws.Range("H86").Value = 2212
$ws.Range("J86").Value = 2201
$ws.Range("K86").Value = 2215.6667
$ws.Range("I89").Value = 2215.6667
$ws.Range("M89").Value = -5462.333500000001
$ws.Range("L89").Value = 11005
$ws.Range("H89").Value = 2212
$ws.Range("J89").Value = 2201
$ws.Range("K89").Value = 11078.3335
$ws.Range("N89").Value = -22237
$ws.Range("K94").Value = 424.5
$ws.Range("I94").Value = 424.5
$ws.Range("M94").Value = 26.5
$ws.Range("H94").Value = 424.5
$ws.Range("H95").Value = 21714.6
$ws.Range("J95").Value = 22857.666
$ws.Range("N95").Value = -28349.666
$ws.Range("L95").Value = 22857.666
$ws.Range("N100").Value = -11717.333
$ws.Range("H100").Value = 5993.385
$ws.Range("L100").Value = 10635.333
$ws.Range("J100").Value = 10635.333
$ws.Range("K107").Value = 1265.5714
$ws.Range("N107").Value = -4236
$ws.Range("I107").Value = 1265.5714
$ws.Range("M107").Value = 654.4286
$ws.Range("H107").Value = 1156.875
$ws.Range("L107").Value = 396
$ws.Range("J107").Value = 396
$ws.Range("N113").Value = -14328
$ws.Range("I113").Value = 3493.5
$ws.Range("M113").Value = -239.5
$ws.Range("H113").Value = 7033.364
$ws.Range("L113").Value = 7820
$ws.Range("J113").Value = 7820
$ws.Range("K113").Value = 3493.5
$ws.Range("K116").Value = 5624.3335
$ws.Range("N116").Value = -12765
$ws.Range("I116").Value = 5624.3335
$ws.Range("M116").Value = -2182.3335
$ws.Range("H116").Value = 5832.875
$ws.Range("L116").Value = 5881
$ws.Range("J116").Value = 5881
$ws.Range("N132").Value = -23053.5005
$ws.Range("I132").Value = 6054.826
$ws.Range("M132").Value = -15634.478
$ws.Range("H132").Value = 6043.0347
$ws.Range("L132").Value = 17993.5005
$ws.Range("J132").Value = 5997.8335
$ws.Range("K132").Value = 18164.478
$ws.Range("I137").Value = 128663.93
$ws.Range("M137").Value = -383441.79
$ws.Range("H137").Value = 128663.93
$ws.Range("K137").Value = 385991.79
$ws.Range("H138").Value = 3029.3108
$ws.Range("J138").Value = 4618.628
$ws.Range("K138").Value = 2474.3226
$ws.Range("N138").Value = -24135.884
$ws.Range("I138").Value = 824.7742
$ws.Range("M138").Value = 2665.6774
$ws.Range("L138").Value = 13855.884
$ws.Range("I141").Value = 5651
$ws.Range("M141").Value = -11773
$ws.Range("H141").Value = 11192.549
$ws.Range("L141").Value = 81373.5
$ws.Range("J141").Value = 27124.5
$ws.Range("K141").Value = 16953
$ws.Range("N141").Value = -91733.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -1913.329
$ws.Range("I32").Value = 2200.329
$ws.Range("H32").Value = 4103.0527
$ws.Range("K32").Value = 2200.329
$ws.Range("I45").Value = 168663.17
$ws.Range("M45").Value = -168286.17
$ws.Range("H45").Value = 96724.17999999999
$ws.Range("K45").Value = 168663.17
$ws.Range("H61").Value = 3601.1428
$ws.Range("L61").Value = 4099.4
$ws.Range("J61").Value = 4099.4
$ws.Range("N61").Value = -4523.4
$ws.Range("K61").Value = 3445.4375
$ws.Range("I61").Value = 3445.4375
$ws.Range("M61").Value = -3233.4375
$ws.Range("H74").Value = 210509.31
$ws.Range("L74").Value = 304170.66
$ws.Range("J74").Value = 304170.66
$ws.Range("K74").Value = 130228.14
$ws.Range("N74").Value = -305918.66
$ws.Range("I74").Value = 130228.14
$ws.Range("M74").Value = -129354.14
$ws.Range("K77").Value = 651140.7
$ws.Range("N77").Value = -1529589.3
$ws.Range("I77").Value = 130228.14
$ws.Range("M77").Value = -646772.7
$ws.Range("H77").Value = 210509.31
$ws.Range("L77").Value = 1520853.3
$ws.Range("J77").Value = 304170.66
$ws.Range("K97").Value = 10990.588
$ws.Range("N97").Value = -22014
$ws.Range("I97").Value = 10990.588
$ws.Range("M97").Value = -10494.588
$ws.Range("H97").Value = 13270.454
$ws.Range("L97").Value = 21022
$ws.Range("J97").Value = 21022
$ws.Range("L102").Value = 6935.7
$ws.Range("H102").Value = 3662.2307
$ws.Range("J102").Value = 6935.7
$ws.Range("N102").Value = -10179.7
$ws.Range("H110").Value = 1540.2
$ws.Range("J110").Value = 1979.6
$ws.Range("K110").Value = 1393.7333
$ws.Range("N110").Value = -6069.6
$ws.Range("I110").Value = 1393.7333
$ws.Range("M110").Value = 651.2666999999999
$ws.Range("L110").Value = 1979.6
$ws.Range("K122").Value = 10716.3
$ws.Range("I122").Value = 3572.1
$ws.Range("M122").Value = -8266.299999999999
$ws.Range("H122").Value = 3611
$ws.Range("N132").Value = -16997
$ws.Range("I132").Value = 2512
$ws.Range("M132").Value = -5006
$ws.Range("H132").Value = 2975.2632
$ws.Range("L132").Value = 11937
$ws.Range("J132").Value = 3979
$ws.Range("K132").Value = 7536
$ws.Range("K136").Value = 10336.3125
$ws.Range("N136").Value = -17398.2
$ws.Range("I136").Value = 3445.4375
$ws.Range("M136").Value = -7786.3125
$ws.Range("H136").Value = 3601.1428
$ws.Range("L136").Value = 12298.2
$ws.Range("J136").Value = 4099.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K28").Value = 5000
$ws.Range("N28").Value = -10585
$ws.Range("M28").Value = -4706
$ws.Range("I28").Value = 5000
$ws.Range("H28").Value = 7498.5
$ws.Range("L28").Value = 9997
$ws.Range("J28").Value = 9997
$ws.Range("N63").Value = -28872
$ws.Range("H63").Value = 27500
$ws.Range("L63").Value = 27500
$ws.Range("J63").Value = 27500
$ws.Range("H66").Value = 27500
$ws.Range("L66").Value = 82500
$ws.Range("J66").Value = 27500
$ws.Range("N66").Value = -89364
$ws.Range("N68").Value = -76759.5
$ws.Range("H68").Value = 75137.5
$ws.Range("L68").Value = 75137.5
$ws.Range("J68").Value = 75137.5
$ws.Range("H71").Value = 75137.5
$ws.Range("L71").Value = 225412.5
$ws.Range("J71").Value = 75137.5
$ws.Range("N71").Value = -233524.5
$ws.Range("N86").Value = -7645.2
$ws.Range("I86").Value = 13332.818
$ws.Range("M86").Value = -12209.818
$ws.Range("L86").Value = 5399.2
$ws.Range("H86").Value = 10853.5625
$ws.Range("J86").Value = 5399.2
$ws.Range("K86").Value = 13332.818
$ws.Range("I89").Value = 13332.818
$ws.Range("M89").Value = -61048.09
$ws.Range("L89").Value = 26996
$ws.Range("H89").Value = 10853.5625
$ws.Range("J89").Value = 5399.2
$ws.Range("K89").Value = 66664.09
$ws.Range("N89").Value = -38228
$ws.Range("K97").Value = 13789.643
$ws.Range("I97").Value = 13789.643
$ws.Range("M97").Value = -12798.643
$ws.Range("H97").Value = 13434.8
$ws.Range("H99").Value = 52182.95
$ws.Range("K99").Value = 64447.062
$ws.Range("I99").Value = 64447.062
$ws.Range("M99").Value = -62949.062
$ws.Range("J101").Value = 0
$ws.Range("N101").Value = $null
$ws.Range("H101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("J105").Value = 2377.6
$ws.Range("K105").Value = 3928.2307
$ws.Range("N105").Value = -5871.6
$ws.Range("I105").Value = 3928.2307
$ws.Range("M105").Value = -2181.2307
$ws.Range("H105").Value = 3497.5
$ws.Range("L105").Value = 2377.6
$ws.Range("H108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("N132").Value = -99117
$ws.Range("H132").Value = 88997
$ws.Range("L132").Value = 88997
$ws.Range("J132").Value = 88997
$ws.Range("K134").Value = 9852.428400000001
$ws.Range("N134").Value = -36902.625
$ws.Range("I134").Value = 3284.1428
$ws.Range("M134").Value = -7317.428400000001
$ws.Range("H134").Value = 7191.7334
$ws.Range("L134").Value = 31832.625
$ws.Range("J134").Value = 10610.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1823.9131
$ws.Range("L16").Value = 2329.5
$ws.Range("J16").Value = 2329.5
$ws.Range("K16").Value = 1435
$ws.Range("N16").Value = -2903.5
$ws.Range("I16").Value = 1435
$ws.Range("M16").Value = -1148
$ws.Range("K31").Value = 2686.5715
$ws.Range("I31").Value = 2686.5715
$ws.Range("M31").Value = -2391.5715
$ws.Range("H31").Value = 25049.262
$ws.Range("K34").Value = 2686.5715
$ws.Range("I34").Value = 2686.5715
$ws.Range("M34").Value = -2484.5715
$ws.Range("H34").Value = 25049.262
$ws.Range("H50").Value = 5860.826
$ws.Range("L50").Value = 5860.826
$ws.Range("J50").Value = 5860.826
$ws.Range("N50").Value = -7110.826
$ws.Range("M58").Value = -1196.8334
$ws.Range("H58").Value = 1932.2894
$ws.Range("L58").Value = 3929
$ws.Range("J58").Value = 3929
$ws.Range("K58").Value = 1399.8334
$ws.Range("N58").Value = -4335
$ws.Range("I58").Value = 1399.8334
$ws.Range("K97").Value = 0
$ws.Range("N97").Value = -61873.5
$ws.Range("I97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("H97").Value = 59891.5
$ws.Range("L97").Value = 59891.5
$ws.Range("J97").Value = 59891.5
$ws.Range("K107").Value = 1324.1936
$ws.Range("N107").Value = -5399.375
$ws.Range("I107").Value = 1324.1936
$ws.Range("M107").Value = 595.8063999999999
$ws.Range("H107").Value = 1372.4359
$ws.Range("L107").Value = 1559.375
$ws.Range("J107").Value = 1559.375
$ws.Range("J109").Value = 59994
$ws.Range("N109").Value = -62074
$ws.Range("H109").Value = 59994
$ws.Range("L109").Value = 59994
$ws.Range("N113").Value = -6669.5
$ws.Range("I113").Value = 1435
$ws.Range("M113").Value = 735
$ws.Range("H113").Value = 1823.9131
$ws.Range("L113").Value = 2329.5
$ws.Range("J113").Value = 2329.5
$ws.Range("K113").Value = 1435
$ws.Range("K122").Value = 3294
$ws.Range("N122").Value = -17062
$ws.Range("I122").Value = 1098
$ws.Range("M122").Value = -844
$ws.Range("H122").Value = 2254.6956
$ws.Range("L122").Value = 12162
$ws.Range("J122").Value = 4054
$ws.Range("N132").Value = -195980
$ws.Range("I132").Value = 1928.7333
$ws.Range("M132").Value = -3256.199900000001
$ws.Range("H132").Value = 32784.367
$ws.Range("L132").Value = 190920
$ws.Range("J132").Value = 63640
$ws.Range("K132").Value = 5786.199900000001
$ws.Range("K134").Value = 237792.258
$ws.Range("N134").Value = -20200.8
$ws.Range("I134").Value = 79264.086
$ws.Range("M134").Value = -235257.258
$ws.Range("H134").Value = 45527.5
$ws.Range("L134").Value = 15130.8
$ws.Range("J134").Value = 5043.6
$ws.Range("K136").Value = 4199.5002
$ws.Range("N136").Value = -16887
$ws.Range("I136").Value = 1399.8334
$ws.Range("M136").Value = -1649.5002
$ws.Range("H136").Value = 1932.2894
$ws.Range("L136").Value = 11787
$ws.Range("J136").Value = 3929

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I7").Value = 90910390
$ws.Range("M7").Value = -272731058
$ws.Range("H7").Value = 66668052
$ws.Range("K7").Value = 272731170
$ws.Range("N12").Value = -555.47059
$ws.Range("H12").Value = 40488.227
$ws.Range("L12").Value = 209.47059
$ws.Range("J12").Value = 69.82353000000001
$ws.Range("H37").Value = 49450
$ws.Range("L37").Value = 148350
$ws.Range("J37").Value = 49450
$ws.Range("N37").Value = -148574
$ws.Range("K46").Value = 5002249.5
$ws.Range("N46").Value = -6405.049999999999
$ws.Range("I46").Value = 1667416.5
$ws.Range("M46").Value = -5002158.5
$ws.Range("H46").Value = 153469.1
$ws.Range("L46").Value = 6223.049999999999
$ws.Range("J46").Value = 2074.35
$ws.Range("K70").Value = 2109
$ws.Range("I70").Value = 703
$ws.Range("M70").Value = -1794
$ws.Range("H70").Value = 703
$ws.Range("K73").Value = 2109
$ws.Range("I73").Value = 703
$ws.Range("M73").Value = -1017
$ws.Range("H73").Value = 703
$ws.Range("K107").Value = 779.6999999999999
$ws.Range("N107").Value = -6245.1429
$ws.Range("I107").Value = 259.9
$ws.Range("M107").Value = 1140.3
$ws.Range("H107").Value = 575.9583
$ws.Range("L107").Value = 2405.1429
$ws.Range("J107").Value = 801.7143
$ws.Range("H131").Value = 7312975
$ws.Range("K131").Value = 10873464
$ws.Range("I131").Value = 3624488
$ws.Range("M131").Value = -10868424
$ws.Range("N132").Value = -35294.60000000001
$ws.Range("I132").Value = 1126.4
$ws.Range("M132").Value = -7607.6
$ws.Range("H132").Value = 2721.4
$ws.Range("L132").Value = 30234.6
$ws.Range("J132").Value = 3359.4
$ws.Range("K132").Value = 10137.6
$ws.Range("N137").Value = -23197.0005
$ws.Range("H137").Value = 2157.5557
$ws.Range("L137").Value = 12997.0005
$ws.Range("J137").Value = 4332.3335
$ws.Range("L140").Value = 17997
$ws.Range("J140").Value = 5999
$ws.Range("K140").Value = 5339.6667
$ws.Range("N140").Value = -28357
$ws.Range("I140").Value = 1779.8889
$ws.Range("M140").Value = -159.6666999999998
$ws.Range("H140").Value = 2201.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M2").Value = 66.85714300000001
$ws.Range("H2").Value = 76923200
$ws.Range("K2").Value = 46.142857
$ws.Range("I2").Value = 46.142857
$ws.Range("N32").Value = -24591
$ws.Range("H32").Value = 23999
$ws.Range("L32").Value = 23999
$ws.Range("J32").Value = 23999
$ws.Range("M58").Value = -23723
$ws.Range("H58").Value = 24749.5
$ws.Range("K58").Value = 24000
$ws.Range("I58").Value = 24000
$ws.Range("K59").Value = 9000
$ws.Range("I59").Value = 9000
$ws.Range("M59").Value = -8417
$ws.Range("H59").Value = 9555
$ws.Range("H80").Value = 31377756
$ws.Range("L80").Value = 402770
$ws.Range("J80").Value = 402770
$ws.Range("K80").Value = 45457296
$ws.Range("M80").Value = -45456298
$ws.Range("I80").Value = 45457296
$ws.Range("N80").Value = -404766
$ws.Range("J83").Value = 402770
$ws.Range("K83").Value = 227286480
$ws.Range("N83").Value = -2023834
$ws.Range("I83").Value = 45457296
$ws.Range("M83").Value = -227281488
$ws.Range("H83").Value = 31377756
$ws.Range("L83").Value = 2013850
$ws.Range("K97").Value = 1614.579
$ws.Range("N97").Value = -2772.2
$ws.Range("I97").Value = 1614.579
$ws.Range("M97").Value = -1118.579
$ws.Range("H97").Value = 1649.0834
$ws.Range("L97").Value = 1780.2
$ws.Range("J97").Value = 1780.2
$ws.Range("K107").Value = 1135.3125
$ws.Range("N107").Value = -4409.5
$ws.Range("I107").Value = 1135.3125
$ws.Range("M107").Value = 784.6875
$ws.Range("H107").Value = 917.6923
$ws.Range("L107").Value = 569.5
$ws.Range("J107").Value = 569.5
$ws.Range("J109").Value = 46071.25
$ws.Range("N109").Value = -48151.25
$ws.Range("H109").Value = 46071.25
$ws.Range("L109").Value = 46071.25
$ws.Range("J120").Value = 37905.75
$ws.Range("N120").Value = -47581.75
$ws.Range("H120").Value = 37905.75
$ws.Range("L120").Value = 37905.75
$ws.Range("J126").Value = 4331.6665
$ws.Range("N126").Value = -17934.9995
$ws.Range("H126").Value = 3538.5715
$ws.Range("L126").Value = 12994.9995
$ws.Range("N132").Value = -14816.75
$ws.Range("I132").Value = 3168.2593
$ws.Range("M132").Value = -6974.777900000001
$ws.Range("H132").Value = 3179.0967
$ws.Range("L132").Value = 9756.75
$ws.Range("J132").Value = 3252.25
$ws.Range("K132").Value = 9504.777900000001
$ws.Range("N133").Value = -130899.336
$ws.Range("H133").Value = 120779.336
$ws.Range("L133").Value = 120779.336
$ws.Range("J133").Value = 120779.336
$ws.Range("N135").Value = -133471.25
$ws.Range("H135").Value = 123331.25
$ws.Range("L135").Value = 123331.25
$ws.Range("J135").Value = 123331.25
$ws.Range("L140").Value = 106900
$ws.Range("J140").Value = 106900
$ws.Range("N140").Value = -117260
$ws.Range("H140").Value = 106900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K22").Value = 128190.29
$ws.Range("I22").Value = 128190.29
$ws.Range("M22").Value = -127895.29
$ws.Range("H22").Value = 128190.29
$ws.Range("K27").Value = 128190.29
$ws.Range("I27").Value = 128190.29
$ws.Range("M27").Value = -128083.29
$ws.Range("H27").Value = 128190.29
$ws.Range("N40").Value = -14034.2
$ws.Range("I40").Value = 5410.76
$ws.Range("M40").Value = -5274.76
$ws.Range("H40").Value = 6802.6665
$ws.Range("L40").Value = 13762.2
$ws.Range("J40").Value = 13762.2
$ws.Range("K40").Value = 5410.76
$ws.Range("K46").Value = 1749.6
$ws.Range("N46").Value = -4757.4546
$ws.Range("I46").Value = 1749.6
$ws.Range("M46").Value = -1561.6
$ws.Range("H46").Value = 3559
$ws.Range("L46").Value = 4381.4546
$ws.Range("J46").Value = 4381.4546
$ws.Range("H55").Value = 1590.35
$ws.Range("L55").Value = 1598.625
$ws.Range("J55").Value = 1598.625
$ws.Range("N55").Value = -1944.625
$ws.Range("I82").Value = 58824580
$ws.Range("M82").Value = -58824219
$ws.Range("H82").Value = 45455572
$ws.Range("L82").Value = 945.4
$ws.Range("J82").Value = 945.4
$ws.Range("K82").Value = 58824580
$ws.Range("N82").Value = -1667.4
$ws.Range("K85").Value = 58824580
$ws.Range("I85").Value = 58824580
$ws.Range("N85").Value = -3441.4
$ws.Range("M85").Value = -58823332
$ws.Range("H85").Value = 45455572
$ws.Range("L85").Value = 945.4
$ws.Range("J85").Value = 945.4
$ws.Range("H93").Value = 1447.85
$ws.Range("L93").Value = 1770.1428
$ws.Range("J93").Value = 1770.1428
$ws.Range("K93").Value = 1274.3077
$ws.Range("M93").Value = -26.30770000000007
$ws.Range("N93").Value = -4266.1428
$ws.Range("I93").Value = 1274.3077
$ws.Range("M100").Value = -3934.8945
$ws.Range("I100").Value = 4475.8945
$ws.Range("H100").Value = 51892.477
$ws.Range("K100").Value = 4475.8945
$ws.Range("J101").Value = 22631.166
$ws.Range("N101").Value = -29121.166
$ws.Range("H101").Value = 22631.166
$ws.Range("L101").Value = 22631.166
$ws.Range("J106").Value = 7079
$ws.Range("K106").Value = 1750
$ws.Range("N106").Value = -9603
$ws.Range("I106").Value = 1750
$ws.Range("M106").Value = -488
$ws.Range("H106").Value = 6190.8335
$ws.Range("L106").Value = 7079
$ws.Range("J109").Value = 39999.5
$ws.Range("N109").Value = -42773.5
$ws.Range("H109").Value = 39999.5
$ws.Range("L109").Value = 39999.5
$ws.Range("N127").Value = -64920
$ws.Range("I127").Value = 0
$ws.Range("M127").Value = $null
$ws.Range("H127").Value = 55000
$ws.Range("L127").Value = 55000
$ws.Range("J127").Value = 55000
$ws.Range("K127").Value = 0
$ws.Range("N132").Value = -26382.8939
$ws.Range("I132").Value = 6647.9756
$ws.Range("M132").Value = -17413.9268
$ws.Range("H132").Value = 6793.533
$ws.Range("L132").Value = 21322.8939
$ws.Range("J132").Value = 7107.6313
$ws.Range("K132").Value = 19943.9268

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K14").Value = 8500
$ws.Range("I14").Value = 8500
$ws.Range("M14").Value = -8332
$ws.Range("H14").Value = 8500
$ws.Range("N22").Value = -596
$ws.Range("H22").Value = 1005
$ws.Range("L22").Value = 10
$ws.Range("J22").Value = 10
$ws.Range("H81").Value = 3526
$ws.Range("K81").Value = 7052
$ws.Range("M81").Value = -5991
$ws.Range("I81").Value = 3526
$ws.Range("I84").Value = 3526
$ws.Range("M84").Value = -29956
$ws.Range("H84").Value = 3526
$ws.Range("K84").Value = 35260
$ws.Range("I96").Value = 4234.8335
$ws.Range("H96").Value = 5723
$ws.Range("K96").Value = 4234.8335
$ws.Range("M96").Value = -2861.8335
$ws.Range("M100").Value = -4059.6666
$ws.Range("I100").Value = 2300.3333
$ws.Range("H100").Value = 2972.75
$ws.Range("K100").Value = 4600.6666
$ws.Range("K107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("H107").Value = 0
$ws.Range("J109").Value = 44777
$ws.Range("N109").Value = -47551
$ws.Range("H109").Value = 44777
$ws.Range("L109").Value = 44777
$ws.Range("I113").Value = 373.5
$ws.Range("M113").Value = 1049.5
$ws.Range("H113").Value = 1443
$ws.Range("K113").Value = 1120.5
$ws.Range("N122").Value = -22237
$ws.Range("H122").Value = 3343.682
$ws.Range("L122").Value = 17337
$ws.Range("J122").Value = 5779
$ws.Range("J126").Value = 3157
$ws.Range("K126").Value = 11273.4999
$ws.Range("N126").Value = -14411
$ws.Range("I126").Value = 3757.8333
$ws.Range("M126").Value = -8803.499899999999
$ws.Range("H126").Value = 3581.1177
$ws.Range("L126").Value = 9471
$ws.Range("N132").Value = -201756.71
$ws.Range("I132").Value = 3566.8262
$ws.Range("M132").Value = -8170.4786
$ws.Range("H132").Value = 18033.2
$ws.Range("L132").Value = 196696.71
$ws.Range("J132").Value = 65565.57000000001
$ws.Range("K132").Value = 10700.4786
$ws.Range("K136").Value = 8562.5625
$ws.Range("N136").Value = -16329
$ws.Range("I136").Value = 2854.1875
$ws.Range("M136").Value = -6012.5625
$ws.Range("H136").Value = 3110.9556
$ws.Range("L136").Value = 11229
$ws.Range("J136").Value = 3743
$ws.Range("H138").Value = 190429
$ws.Range("J138").Value = 190429
$ws.Range("N138").Value = -200709
$ws.Range("L138").Value = 190429
$ws.Range("N139").Value = -87328.336
$ws.Range("L139").Value = 77048.336
$ws.Range("H139").Value = 77048.336
$ws.Range("J139").Value = 77048.336
